$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Version table: "Número de versión" -> "1.0"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Número de versión", $true, $false, $false, $false, $false, $true, 1, $false, "1.0", 2)

# ---------------------------------------------------------------------
# 2. Delivery-date field: "dd/mm/aaaa" -> "14/03/2014"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("dd/mm/aaaa", $true, $false, $false, $false, $false, $true, 1, $false, "14/03/2014", 2)

# ---------------------------------------------------------------------
# 3. Case-of-study field: "Identificación del (sub)caso de estudio" -> "Gestion de tarjetas universitarias"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Identificación del (sub)caso de estudio", $true, $false, $false, $false, $false, $true, 1, $false, "Gestion de tarjetas universitarias", 2)

# ---------------------------------------------------------------------
# 4. Delivery/iteration field: "Identificación de entrega o iteración" -> "1"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Identificación de entrega o iteración", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# ---------------------------------------------------------------------
# 5. Place a fresh "_GoBack" bookmark at the start of "Fecha de entrega"
#    (mirrors Word's automatic last-edit bookmark; the old one near the
#    end of the document is replaced/removed and every ToC bookmark ID
#    shifts up by one).
# ---------------------------------------------------------------------
$goback = $d.Content
$goback.Find.Execute("Fecha de entrega")
$gobackPoint = $d.Range($goback.Start, $goback.Start)
$d.Bookmarks.Add("_GoBack", $gobackPoint)
